$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to remain text for the cells we touch that look like plain numbers,
# by temporarily setting the whole Price column to Text format, then restoring Normal style
# so the saved style index matches the original (no style attribute).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '58.262.75'
$ws.Range('E2').Value = '  -0.13%  '
$ws.Range('D3').Value = '2.595.33'
$ws.Range('E3').Value = '  -0.77%  '
$ws.Range('E4').Value = '  +0.10%  '
$ws.Range('D5').Value = '522.79'
$ws.Range('D6').Value = '143.43'
$ws.Range('E6').Value = '  +0.87%  '
$ws.Range('E7').Value = '  -0.17%  '
$ws.Range('D8').Value = '0.569'
$ws.Range('E8').Value = '  +0.51%  '
$ws.Range('D9').Value = '2.617.26'
$ws.Range('E10').Value = '  -1.49%  '
$ws.Range('D11').Value = '0.101'
$ws.Range('E11').Value = '  -1.30%  '
$ws.Range('D12').Value = '0.343'
$ws.Range('E12').Value = '  +1.92%  '
$ws.Range('E13').Value = '  -0.05%  '
$ws.Range('D14').Value = '3.057.54'
$ws.Range('E14').Value = '  -0.55%  '
$ws.Range('D15').Value = '58.243.29'
$ws.Range('E15').Value = '  -0.11%  '
$ws.Range('D16').Value = '20.40'
$ws.Range('E16').Value = '  -2.39%  '
$ws.Range('E17').Value = '  -1.04%  '
$ws.Range('D18').Value = '2.572.14'
$ws.Range('E18').Value = '  -2.47%  '
$ws.Range('D19').Value = '339.97'
$ws.Range('E19').Value = '  +0.84%  '
$ws.Range('D20').Value = '4.35'
$ws.Range('E20').Value = '  -0.90%  '
$ws.Range('D21').Value = '10.24'
$ws.Range('E21').Value = '  -1.51%  '
$ws.Range('D22').Value = '6.43'
$ws.Range('E22').Value = '  +2.31%  '
$ws.Range('D23').Value = '0.998'
$ws.Range('E23').Value = '  -0.08%  '
$ws.Range('D24').Value = '65.35'
$ws.Range('E24').Value = '  +1.20%  '
$ws.Range('E25').Value = '  +1.10%  '
$ws.Range('E26').Value = '  -2.31%  '
$ws.Range('B27').Value = 'WrappedeETH'
$ws.Range('C27').Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range('D27').Value = '2.718.23'
$ws.Range('E27').Value = '  -0.61%  '
$ws.Range('B28').Value = 'Binance-PegBSC-USD'
$ws.Range('C28').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D28').Value = '0.996'
$ws.Range('E28').Value = '  -0.30%  '
$ws.Range('D29').Value = '7.02'
$ws.Range('E29').Value = '  -1.24%  '
$ws.Range('D30').Value = '0.0₃0747'
$ws.Range('E30').Value = '  -5.25%  '
$ws.Range('D32').Value = '6.18'
$ws.Range('E32').Value = '  -4.61%  '
$ws.Range('D33').Value = '1.58'
$ws.Range('E33').Value = '  -0.40%  '
$ws.Range('D34').Value = '18.79'
$ws.Range('E34').Value = '  +0.24%  '
$ws.Range('D35').Value = '149.80'
$ws.Range('D36').Value = '4.01'
$ws.Range('E36').Value = '  -1.67%  '
$ws.Range('D37').Value = '1.13'
$ws.Range('E37').Value = '  -3.84%  '
$ws.Range('D38').Value = '0.869'
$ws.Range('E38').Value = '  -1.09%  '
$ws.Range('D39').Value = '0.868'
$ws.Range('E39').Value = '  +2.01%  '
$ws.Range('D40').Value = '1.46'
$ws.Range('E40').Value = '  +2.70%  '
$ws.Range('D41').Value = '35.99'
$ws.Range('E41').Value = '  -0.67%  '
$ws.Range('E42').Value = '  -2.03%  '
$ws.Range('E43').Value = '  -0.24%  '
$ws.Range('D44').Value = '0.604'
$ws.Range('E44').Value = '  +0.21%  '
$ws.Range('D45').Value = '270.90'
$ws.Range('E45').Value = '  +0.92%  '
$ws.Range('D46').Value = '0.0956'
$ws.Range('E46').Value = '  -0.98%  '
$ws.Range('D47').Value = '10.67'
$ws.Range('E47').Value = '  +0.24%  '
$ws.Range('D48').Value = '18.78'
$ws.Range('E48').Value = '  -1.63%  '
$ws.Range('D49').Value = '0.0523'
$ws.Range('E49').Value = '  -1.43%  '
$ws.Range('B50').Value = 'InjectiveProtocol'
$ws.Range('C50').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D50').Value = '18.78'
$ws.Range('E50').Value = '  +3.39%  '
$ws.Range('D51').Value = '1.970.61'
$ws.Range('E51').Value = '  -2.49%  '

# Restore the original (default/Normal) style for the Price column so no stray
# number-format style lingers on any cell in that column.
$ws.Range("D2:D51").Style = "Normal"
